$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to text first so numeric-looking strings (e.g. "1.01") are not
# auto-converted to numbers by Excel; style is reset back to Normal afterwards so
# the cells keep their original (unstyled) appearance.
$colD = $ws.Range("D2:D51")
$colD.NumberFormat = "@"

$ws.Cells.Item(2, 4).Value = "36.993.15"
$ws.Cells.Item(2, 5).Value = "  -1.46%  "

$ws.Cells.Item(3, 4).Value = "2.028.42"
$ws.Cells.Item(3, 5).Value = "  -2.60%  "

$ws.Cells.Item(4, 4).Value = "1.01"
$ws.Cells.Item(4, 5).Value = "  +0.75%  "

$ws.Cells.Item(5, 4).Value = "227.16"
$ws.Cells.Item(5, 5).Value = "  -2.49%  "

$ws.Cells.Item(6, 5).Value = "  -3.77%  "

$ws.Cells.Item(7, 5).Value = "  +0.06%  "

$ws.Cells.Item(8, 4).Value = "55.13"
$ws.Cells.Item(8, 5).Value = "  -4.81%  "

$ws.Cells.Item(9, 4).Value = "0.382"
$ws.Cells.Item(9, 5).Value = "  -2.26%  "

$ws.Cells.Item(10, 4).Value = "0.0794"
$ws.Cells.Item(10, 5).Value = "  +1.74%  "

$ws.Cells.Item(11, 5).Value = "  -3.59%  "

$ws.Cells.Item(12, 4).Value = "2.325.39"
$ws.Cells.Item(12, 5).Value = "  -2.71%  "

$ws.Cells.Item(13, 4).Value = "14.36"
$ws.Cells.Item(13, 5).Value = "  -4.68%  "

$ws.Cells.Item(14, 4).Value = "20.57"
$ws.Cells.Item(14, 5).Value = "  -2.54%  "

$ws.Cells.Item(15, 4).Value = "0.749"
$ws.Cells.Item(15, 5).Value = "  -3.05%  "

$ws.Cells.Item(16, 4).Value = "5.17"
$ws.Cells.Item(16, 5).Value = "  -3.64%  "

$ws.Cells.Item(17, 4).Value = "2.024.10"
$ws.Cells.Item(17, 5).Value = "  -2.67%  "

$ws.Cells.Item(18, 4).Value = "36.973.09"
$ws.Cells.Item(18, 5).Value = "  -1.49%  "

$ws.Cells.Item(19, 4).Value = "6.16"
$ws.Cells.Item(19, 5).Value = "  +1.82%  "

$ws.Cells.Item(20, 4).Value = "68.92"
$ws.Cells.Item(20, 5).Value = "  -2.57%  "

$ws.Cells.Item(21, 4).Value = "0.0₃0826"
$ws.Cells.Item(21, 5).Value = "  -0.85%  "

$ws.Cells.Item(22, 5).Value = "  -0.97%  "

$ws.Cells.Item(23, 4).Value = "1.00"
$ws.Cells.Item(23, 5).Value = "  +0.06%  "

$ws.Cells.Item(24, 5).Value = "  +3.70%  "

$ws.Cells.Item(25, 4).Value = "2.27"
$ws.Cells.Item(25, 5).Value = "  -4.73%  "

$ws.Cells.Item(26, 4).Value = "166.66"
$ws.Cells.Item(26, 5).Value = "  -2.17%  "

$ws.Cells.Item(27, 4).Value = "9.25"
$ws.Cells.Item(27, 5).Value = "  -4.55%  "

$ws.Cells.Item(28, 4).Value = "0.127"
$ws.Cells.Item(28, 5).Value = "  -4.88%  "

$ws.Cells.Item(29, 4).Value = "18.78"
$ws.Cells.Item(29, 5).Value = "  -4.62%  "

$ws.Cells.Item(30, 4).Value = "1.34"
$ws.Cells.Item(30, 5).Value = "  -3.41%  "

$ws.Cells.Item(31, 5).Value = "  -4.77%  "

$ws.Cells.Item(32, 4).Value = "4.49"
$ws.Cells.Item(32, 5).Value = "  -3.55%  "

$ws.Cells.Item(33, 4).Value = "0.0616"
$ws.Cells.Item(33, 5).Value = "  -3.19%  "

$ws.Cells.Item(34, 4).Value = "4.45"
$ws.Cells.Item(34, 5).Value = "  -4.81%  "

$ws.Cells.Item(35, 4).Value = "2.37"
$ws.Cells.Item(35, 5).Value = "  -4.45%  "

$ws.Cells.Item(36, 4).Value = "1.84"
$ws.Cells.Item(36, 5).Value = "  +1.12%  "

$ws.Cells.Item(37, 4).Value = "1.01"
$ws.Cells.Item(37, 5).Value = "  +0.48%  "

$ws.Cells.Item(38, 5).Value = "  -4.38%  "

$ws.Cells.Item(39, 4).Value = "5.35"
$ws.Cells.Item(39, 5).Value = "  +0.77%  "

$ws.Cells.Item(40, 5).Value = "  -5.18%  "

$ws.Cells.Item(41, 2).Value = "Maker"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Cells.Item(41, 4).Value = "1.492.76"
$ws.Cells.Item(41, 5).Value = "  +1.95%  "

$ws.Cells.Item(42, 2).Value = "InjectiveProtocol"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Cells.Item(42, 4).Value = "17.07"
$ws.Cells.Item(42, 5).Value = "  +1.61%  "

$ws.Cells.Item(43, 4).Value = "95.38"
$ws.Cells.Item(43, 5).Value = "  -4.98%  "

$ws.Cells.Item(44, 4).Value = "0.0931"
$ws.Cells.Item(44, 5).Value = "  -2.70%  "

$ws.Cells.Item(45, 2).Value = "TrustWalletToken"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(45, 4).Value = "1.14"
$ws.Cells.Item(45, 5).Value = "  -5.22%  "

$ws.Cells.Item(46, 2).Value = "HuobiToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(46, 4).Value = "2.75"
$ws.Cells.Item(46, 5).Value = "  -5.63%  "

$ws.Cells.Item(48, 4).Value = "1.01"
$ws.Cells.Item(48, 5).Value = "  -4.20%  "

$ws.Cells.Item(50, 4).Value = "2.212.45"
$ws.Cells.Item(50, 5).Value = "  -2.67%  "

$ws.Cells.Item(51, 4).Value = "3.59"
$ws.Cells.Item(51, 5).Value = "  -8.57%  "

# Restore default (Normal) style on column D so no stray number-format styling
# remains attached to the text cells.
$colD.Style = "Normal"

